# week10Internet.pptx — slide 6 "Coding Challenge":
#  - grow the body placeholder's height
#  - "...in the string..." -> "...in an array..."
#  - "aaa " (as two runs) -> single run "[a, a, a] "
#  - "aab " (as two runs) -> single run "[a, a, b] "
#  - "abba" -> "[a, b, b, a]"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shape = $s.Shapes.Item(2)

# --- 1. Resize the placeholder (cy: 5942278 -> 6719981 EMU). Shape.Height is
#        expressed in points (1 pt = 12700 EMU) over this COM surface.
$shape.Height = 6719981 / 12700

$tf = $shape.TextFrame
$tr = $tf.TextRange

# --- 2. "Write a loop ... in the string. For example" -> "...in an array..."
#        This whole sentence lives in a single run; rewrite that run in one
#        shot with a literal replacement string so it stays a single run
#        afterwards. NOTE: the run contains curly quotes (U+2018/U+2019)
#        around "number of changes" — the COM text marshalling here folds
#        non-ASCII characters to ASCII look-alikes on *read*, so the
#        replacement text is spelled out literally (with the proper
#        Unicode quote characters) instead of being derived from a
#        `.Text` value that was read out of the run.
$full = $tr.Text
$startIdx = $full.IndexOf("Write a loop that counts")
$endIdx = $full.IndexOf("aaa") - 1   # back up over the paragraph-break (chr 13)
$sentence = $tr.Characters($startIdx + 1, $endIdx - $startIdx)
$lq = [char]0x2018
$rq = [char]0x2019
$sentence.Text = "Write a loop that counts the ${lq}number of changes${rq} in an array. For example"

# --- 3. "aaa" (err-flagged run) + " " (plain run) -> one plain run "[a, a, a] "
#        Write the full replacement into the *second* (non err) run first,
#        then blank out the first (err-flagged) run so only clean formatting
#        survives and the two runs collapse into one.
$full = $tr.Text
$idx = $full.IndexOf("aaa")
$spaceRun = $tr.Characters($idx + 3 + 1, 1)
$spaceRun.Text = "[a, a, a] "
$wordRun = $tr.Characters($idx + 1, 3)
$wordRun.Text = ""

# --- 4. "aab" (err-flagged run) + " " (plain run) -> one plain run "[a, a, b] "
$full = $tr.Text
$idx = $full.IndexOf("aab")
$spaceRun = $tr.Characters($idx + 3 + 1, 1)
$spaceRun.Text = "[a, a, b] "
$wordRun = $tr.Characters($idx + 1, 3)
$wordRun.Text = ""

# --- 5. "abba" -> "[a, b, b, a]"
$full = $tr.Text
$idx = $full.IndexOf("abba")
$wordRun = $tr.Characters($idx + 1, 4)
$wordRun.Text = "[a, b, b, a]"
